# Edit: renumber idMedico (Medico sheet) starting at 1 instead of 2, and
# propagate the new idMedico values to the idMedico foreign-key column on
# the Consulta sheet. Also update the active sheet/selection state to match
# (Consulta becomes the active/selected tab, with a new selection), and the
# Medico sheet's own stored selection moves to A4.

$wb = $excel.ActiveWorkbook

# --- Medico sheet: renumber idMedico column (A) ---
$wsMedico = $wb.Worksheets.Item("Medico")
$wsMedico.Range("A2").Value = 1
$wsMedico.Range("A3").Value = 2
$wsMedico.Range("A4").Value = 3

# --- Consulta sheet: update idMedico foreign-key column (C) to match ---
$wsConsulta = $wb.Worksheets.Item("Consulta")
$wsConsulta.Range("C2").Value = 3
$wsConsulta.Range("C3").Value = 2
$wsConsulta.Range("C4").Value = 2
$wsConsulta.Range("C5").Value = 2
$wsConsulta.Range("C6").Value = 1
$wsConsulta.Range("C7").Value = 3
$wsConsulta.Range("C8").Value = 1

# --- Selection bookkeeping ---
# Medico sheet's saved selection moves to A4 (no change of active sheet).
$wsMedico.Range("A4").Select()

# Consulta becomes the active/selected sheet, with its selection at G23.
$wsConsulta.Activate()
$wsConsulta.Range("G23").Select()
